# "Error Calculations and Plots"
#
# The missing-data sample sheet was regenerated: the "RM 232" row and the
# "SC 92" row are dropped entirely (shrinking the sheet from A1:F35 to
# A1:F33), and three cells flip between "missing" and "imputed" for the
# rows that remain:
#   - SC 5   column A (col B of the table) goes from missing -> -20.2
#   - SC 101 column A (col B of the table) goes from -20.4   -> missing
#   - SC 232 column C (col D of the table) goes from missing -> -14.1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "RM 232" row (row 26) completely; everything below shifts up.
$ws.Rows("26").Delete()

# Remove the "SC 92" row, which is now row 27 after the shift above.
$ws.Rows("27").Delete()

# Fill in the newly-imputed value for "SC 5" (now row 26).
$ws.Range("B26").Value = -20.2

# "SC 101" (now row 27) becomes missing in column B.
$ws.Range("B27").ClearContents()

# Fill in the newly-imputed value for "SC 232" (now row 33), column D.
$ws.Range("D33").Value = -14.1
